$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AC / AD / AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header style (bold, bordered, centered) used by the existing header row
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null

# Team record totals for every data row (2 through 38)
for ($row = 2; $row -le 38; $row++) {
    $ws.Cells.Item($row, 29).Value = 53
    $ws.Cells.Item($row, 30).Value = 64
    $ws.Cells.Item($row, 31).Value = 0
}
